# Apply cell updates per commit: "Updated cryptos list on Fri Sep 22 22:25:57 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells that look like plain numbers need an explicit text format first, otherwise
# Excel would silently convert the assigned string into a Number and drop formatting such
# as trailing zeros (e.g. "0.910" -> 0.91). Every Price/Volume cell in this sheet is stored
# as literal text, so we force that per cell before writing the new value.

$ws.Range("D2").Value = '26.620.64'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").Value = '1.596.61'
$ws.Range("E3").Value = '  +0.50%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.62'
$ws.Range("E5").Value = '  +0.26%  '
$ws.Range("E6").Value = '  +1.13%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  +0.31%  '
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.49'
$ws.Range("E10").Value = '  -0.56%  '
$ws.Range("E11").Value = '  +0.36%  '
$ws.Range("D12").Value = '1.820.56'
$ws.Range("E12").Value = '  +0.52%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.03'
$ws.Range("E13").Value = '  +0.17%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.541.11'
$ws.Range("E14").Value = '  -3.04%  '
$ws.Range("E15").Value = '  -0.02%  '
$ws.Range("E16").Value = '  -0.30%  '
$ws.Range("D17").Value = '26.609.29'
$ws.Range("E17").Value = '  -0.02%  '
$ws.Range("D18").Value = '0.0₃0732'
$ws.Range("E18").Value = '  +0.67%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '208.73'
$ws.Range("E19").Value = '  +0.22%  '
$ws.Range("E20").Value = '  +0.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.95'
$ws.Range("E21").Value = '  +3.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.27'
$ws.Range("E22").Value = '  +0.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.30'
$ws.Range("E23").Value = '  -1.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.91'
$ws.Range("E24").Value = '  +0.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.45'
$ws.Range("E25").Value = '  -0.84%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.14'
$ws.Range("E27").Value = '  -1.37%  '
$ws.Range("E28").Value = '  +0.80%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.27'
$ws.Range("E29").Value = '  -0.15%  '
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("E31").Value = '  +0.58%  '
$ws.Range("E32").Value = '  +0.26%  '
$ws.Range("E33").Value = '  -1.72%  '
$ws.Range("E34").Value = '  +1.08%  '
$ws.Range("D35").Value = '1.282.21'
$ws.Range("E35").Value = '  -1.95%  '
$ws.Range("E36").Value = '  +1.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.50'
$ws.Range("E37").Value = '  +1.04%  '
$ws.Range("E38").Value = '  -0.25%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.844'
$ws.Range("E39").Value = '  +1.85%  '
$ws.Range("E40").Value = '  +0.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.47'
$ws.Range("E41").Value = '  +1.94%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.786'
$ws.Range("E42").Value = '  -0.76%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '64.40'
$ws.Range("E43").Value = '  +2.81%  '
$ws.Range("B44").Value = 'MXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.19'
$ws.Range("E44").Value = '  +1.12%  '
$ws.Range("D45").Value = '1.733.11'
$ws.Range("E45").Value = '  +0.51%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.910'
$ws.Range("E46").Value = '  +8.79%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '89.71'
$ws.Range("E47").Value = '  +0.11%  '
$ws.Range("E48").Value = '  -0.22%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.103'
$ws.Range("E49").Value = '  +5.09%  '
$ws.Range("E50").Value = '  +0.45%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.49'
$ws.Range("E51").Value = '  +0.09%  '
